$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Logs")

# --- Add new row 5: "Opvolging retour" ---
$ws.Range("A5").Value = "Opvolging retour"
$ws.Range("B5").Value = "mailmind.test@zohomail.eu"
$ws.Range("D5").Value = "Retour / Terugbetaling"
$ws.Range("F5").Value = "2025-08-28 18:04:08"
$ws.Range("G5").Value = "Nee"
$ws.Range("H5").Value = "Ja"
$ws.Range("I5").Value = "Nee"
$ws.Range("J5").Value = "Nee"

# --- Add new row 6: "Retour status" ---
$ws.Range("A6").Value = "Retour status"
$ws.Range("B6").Value = "mailmind.test@zohomail.eu"
$ws.Range("D6").Value = "Retour / Terugbetaling"
$ws.Range("F6").Value = "2025-08-28 18:04:08"
$ws.Range("G6").Value = "Ja"
$ws.Range("H6").Value = "Nee"
$ws.Range("I6").Value = "Nee"
$ws.Range("J6").Value = "Nee"

# --- Extend conditional formatting ranges to include the new rows ---
$fcsD = $ws.Range("D2:D4").FormatConditions
for ($i = 1; $i -le $fcsD.Count; $i++) {
    $fcsD.Item($i).ModifyAppliesToRange($ws.Range("D2:D6"))
}

$fcsG = $ws.Range("G2:G4").FormatConditions
for ($i = 1; $i -le $fcsG.Count; $i++) {
    $fcsG.Item($i).ModifyAppliesToRange($ws.Range("G2:G6"))
}

$fcsH = $ws.Range("H2:H4").FormatConditions
for ($i = 1; $i -le $fcsH.Count; $i++) {
    $fcsH.Item($i).ModifyAppliesToRange($ws.Range("H2:H6"))
}

$fcsI = $ws.Range("I2:I4").FormatConditions
for ($i = 1; $i -le $fcsI.Count; $i++) {
    $fcsI.Item($i).ModifyAppliesToRange($ws.Range("I2:I6"))
}

$fcsJ = $ws.Range("J2:J4").FormatConditions
for ($i = 1; $i -le $fcsJ.Count; $i++) {
    $fcsJ.Item($i).ModifyAppliesToRange($ws.Range("J2:J6"))
}

# --- Update Dashboard summary count ---
$dash = $wb.Worksheets.Item("Dashboard")
$dash.Range("B2").Value = 5
